$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values ---------------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Formatting (bold, centered/top, thin box border) on B1 and A2 --------
$cell1 = $ws.Range("B1")
$cell1.Font.Bold = $true
$cell1.HorizontalAlignment = -4108   # xlCenter
$cell1.VerticalAlignment = -4160     # xlTop
$cell1.Borders.LineStyle = 1         # xlContinuous
$cell1.Borders.Weight = 2            # xlThin

# Copy the resulting format onto A2 so both cells resolve to the exact same
# cell-style (xf) entry instead of each independently accruing their own.
$cell1.Copy()
$cell2 = $ws.Range("A2")
$cell2.PasteSpecial(-4122)           # xlPasteFormats

$excel.CutCopyMode = $false
